$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# Cells are formatted as Text first so numeric-looking strings (e.g. "259.48")
# are stored as literal text (matching the source data), then the temporary
# Text number-format is cleared back to the default "Normal" style so no
# lingering style index is left on the cell.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D29","D30","D31","D32","D33","D34","D35","D36","D38","D39","D40","D41","D42","D43","D44","D45","D47","D48","D49","D51")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.628.10"
$ws.Range("D3").Value = "1.843.24"
$ws.Range("D4").Value = "1.000"
$ws.Range("D5").Value = "259.48"
$ws.Range("D6").Value = "1.000"
$ws.Range("D7").Value = "0.5271"
$ws.Range("D8").Value = "0.3148"
$ws.Range("D9").Value = "0.06808"
$ws.Range("D11").Value = "0.7818"
$ws.Range("D12").Value = "0.07759"
$ws.Range("D13").Value = "1.836.75"
$ws.Range("D14").Value = "88.03"
$ws.Range("D15").Value = "5.011"
$ws.Range("D16").Value = "0.9999"
$ws.Range("D17").Value = "13.85"
$ws.Range("D18").Value = "0.9999"
$ws.Range("D19").Value = "0.000007934"
$ws.Range("D20").Value = "26.647.29"
$ws.Range("D21").Value = "2.072.36"
$ws.Range("D22").Value = "4.609"
$ws.Range("D23").Value = "5.979"
$ws.Range("D24").Value = "9.321"
$ws.Range("D25").Value = "142.59"
$ws.Range("D26").Value = "2.205"
$ws.Range("D29").Value = "110.91"
$ws.Range("D30").Value = "4.184"
$ws.Range("D31").Value = "0.08732"
$ws.Range("D32").Value = "4.076"
$ws.Range("D33").Value = "0.04876"
$ws.Range("D34").Value = "0.7312"
$ws.Range("D35").Value = "1.141"
$ws.Range("D36").Value = "2.860"
$ws.Range("D38").Value = "2.258"
$ws.Range("D39").Value = "0.01731"
$ws.Range("D40").Value = "0.4795"
$ws.Range("D41").Value = "0.8945"
$ws.Range("D42").Value = "110.03"
$ws.Range("D43").Value = "5.925"
$ws.Range("D44").Value = "1.000"
$ws.Range("D45").Value = "7.667"
$ws.Range("D47").Value = "9.084"
$ws.Range("D48").Value = "0.1240"
$ws.Range("D49").Value = "0.05815"
$ws.Range("D51").Value = "0.8931"

foreach ($ref in $priceCells) {
    $ws.Range($ref).Style = "Normal"
}

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("E8").Value = "  -3.47%  "
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("E34").Value = "  +2.12%  "
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("E51").Value = "  +0.80%  "
